# Fixed d,v,w ion generation
#
# Edits applied:
#  1. losses sheet (B21/B32/B34): "partial_sidechain" -> "partial_sidechain_prime"
#  2. ion_type_list sheet: row13 becomes "d'" (N-based), row14 becomes "v" (was "w"),
#     two new rows are inserted (15 = "w", 16 = "w'"), and the old row15 ("b+H2O")
#     shifts down to row17 unmodified.
#  3. View/selection state: ion_type_list becomes the active tab, amino_acids no
#     longer is; a couple of column widths and pane/selection positions change.

$wb = $excel.ActiveWorkbook

# --- losses sheet: rename the "partial_sidechain" loss label -------------
$wsLosses = $wb.Worksheets.Item("losses")
$wsLosses.Range("B21").Value = "partial_sidechain_prime"
$wsLosses.Range("B32").Value = "partial_sidechain_prime"
$wsLosses.Range("B34").Value = "partial_sidechain_prime"

# losses sheet view/column tweaks
$wsLosses.Columns.Item(2).ColumnWidth = 30.85546875
$panes = $wsLosses.Panes
$wsLosses.Application.ActiveWindow.ScrollRow = 8
$wsLosses.Range("B21").Select()

# --- ion_type_list sheet: fix up d / v / w / w' rows ----------------------
$ws = $wb.Worksheets.Item("ion_type_list")

# row13 was "v" (C-based) -> becomes "d'" (N-based)
$ws.Range("A13").Value = "d'"
$ws.Range("B13").Value = "N"
$ws.Range("C13").Value = -1
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = -1

# row14 was "w" (C-based) -> becomes "v" (C-based)
$ws.Range("A14").Value = "v"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0

# insert two fresh rows at 15 (old row15 "b+H2O" is pushed down to row17)
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()

# new row15: "w" (C-based), same composition the old row14 "w" used to have
$ws.Range("A15").Value = "w"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = -1
$ws.Range("E15").Value = -1
$ws.Range("F15").Value = 0
$ws.Range("G15:T15").Value = 0

# new row16: "w'" (C-based)
$ws.Range("A16").Value = "w'"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = -1
$ws.Range("E16").Value = -1
$ws.Range("F16").Value = 0
$ws.Range("G16:T16").Value = 0

# ion_type_list becomes the active tab/selection
$ws.Range("E19").Select()
$ws.Activate()

# amino_acids sheet: column width split + no longer the active tab
$wsAmino = $wb.Worksheets.Item("amino_acids")
$wsAmino.Columns.Item(5).ColumnWidth = 17.140625
$wsAmino.Columns.Item(6).ColumnWidth = 21.7109375
